$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Alpura Selecta Leche Entera Deslactosada 100% de Vaca, Ultrapasteurizada en paquete de 12 Piezas de 1 Litro cada una, de fácil digestión, Fuente de proteína y Calcio, adicionada con vitaminas A y D. Cremosidad y Sabor superior',
    'Lámpara Escritorio LED,Lámpara de Oficina,Lámpara de Mesa con Atenuación de Protección para Los Ojos,Luz de Lectura con Usb Recargable,3 Colores con Clip,360 ° Flexión,para Trabajar, Estudiar, Leer',
    'Lámpara Escritorio LED,Lámpara de Oficina,Lámpara de Mesa con Atenuación de Protección para Los Ojos,Luz de Lectura con Usb Recargable,3 Colores con Clip,360 ° Flexión,para Trabajar, Estudiar, Leer',
    'Govee Tira LED 20 Metros, Cinta LED WiFi Compatible con Alexa, Sync de Música, Bricolaje de Múltiples Colores en Una Línea para Habitación, Decoracion de Año Nuevo',
    'UGREEN 2M Cable USB C a USB A 2.0, 3A USB C Cable Tipo C Carga Rápida con Trenza de Nailon, Cable USB Tipo C Compatible con iPhone 16 15 Pro Max Plus, Galaxy S24 S23 S22 S21, Redmi, Pixel, Huawei',
    'UGREEN 2M Cable USB C a USB A 2.0, 3A USB C Cable Tipo C Carga Rápida con Trenza de Nailon, Cable USB Tipo C Compatible con iPhone 16 15 Pro Max Plus, Galaxy S24 S23 S22 S21, Redmi, Pixel, Huawei',
    'PUMA - Sudadera de forro polar con cuello redondo para mujer (negro, L), Negro -, Large',
    'Cómo mandar a la mierda de forma educada',
    'Zombieland: Double Tap - Roadtrip - PS4 - Estándar Edition',
    'Milky Way chocolate 14 barras de 48 g c/u, Total 672g',
    'Simple Wisdom Escritorio Minimalista,Escritorio con Repisa En Casa 1.2m,DiseñO Minimalista Moderno Escritorio Oficina,（120cm*50cm） (Blanco Puro)',
    '2pcs Lámpara de Techo 36W ，Wodoso Lampara de Techo LED Empotrable，3000K-6500K，30CM,Con Mando a Distancia, Compatible con app',
    'Porta de Pasaporte de Cuero PU, Funda para Pasaporte con Bloqueo RFID, Billetera de pasaporte, Titular de Pasaporte de Billetera de Viajes con Soporte para la Pluma para Mujeres/Hombres (Amarillo)',
    'Santa Clara - Leche Deslactosada Light, 12 Pack de 1 Lt cada uno',
    'WERA - 05051493001-838 RA-R M Destornillador de bitholding con función de carraca, 1/4", 1/4" x 123,5 mm',
    'SlowTon Quitapelos y Rodillos para Mascotas, Cepillo de Limpieza Removedor de Pelaje para Perros y Gatos, Reutilizable Rodillo Quita Pelos, Eliminador de Pelo para Animales Cepillos para Perros(Azure)',
    'Casio Serie A158WA | Reloj digital unisex | clásico | WR | Cronómetro 100 SEC | Alarma diaria | Mantenimiento de tiempo regular: hora, minuto, segundo, pm, fecha, día | Luz LED | Batería de 7 años',
    'amiibo Toad Super Mario Series',
    'Perras de reserva',
    'ACOVARK Funda para Samsung Galaxy S24 Ultra 5G,con 3piezas de película Protectora de Alta adherencia, magnético Anti-Amarillo Antigolpes Protectora Carcasa,Ultrafino y tecnológico,protección Ocular',
    'Fresh Step Multi-Cat con el Poder de Febreze, Arena para Gatos aglutinante de 14 LB/ 6.35 Kg'
)

$startRow = 41
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}
